# Generate Report for Handback
# This script applies the "handback" localization-status update:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this text is shared by the Overview sheet's E/F columns and the "Status" column
#    on the zh-cn / de-de sheets)
#  - The zh-cn and de-de sheets get their "Latest Target File" (I) and
#    "Latest Handback File" (J) columns populated for both data rows, with I also
#    becoming a hyperlink (same target as the corresponding row's source-file hyperlink
#    in column A).
#  - The "Latest Handback DateTime" (K) column is stamped with the handback timestamp
#    (a different timestamp per language sheet).
#  - A few columns are widened to better fit the newly-populated content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: update status text for both rows (columns E and F) ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn / de-de sheets: update the "Status" column (C) for both rows ---
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# Widen columns E and F on the Overview sheet (closest reachable width to the
# target 29.9777047293527, given this engine's column-width quantization).
$wsOverview.Range("E1").ColumnWidth = 29.2
$wsOverview.Range("F1").ColumnWidth = 29.2

# Populates the "Latest Target File" (I), "Latest Handback File" (J) and
# "Latest Handback DateTime" (K) cells for one data row of a language sheet.
# (Uses only positional parameters - named parameter binding to custom
# functions is not supported by this host.)
function Set-HandbackRow($ws, $row, $addr, $mdName, $xlfName, $dt) {
    # Column I ("Latest Target File"): hyperlink to the same markdown file as the
    # source-file hyperlink in column A, with the md file name as display text.
    $ws.Hyperlinks.Add($ws.Range("I" + $row), $addr, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

    # Column J ("Latest Handback File"): the generated xliff file name.
    $ws.Range("J" + $row).Value = $xlfName

    # Column K ("Latest Handback DateTime"): handback timestamp.
    $ws.Range("K" + $row).Value = $dt
}

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8070046300b1d3e10b448207456cdc9874c6e560/e2e"

$addrC1 = $githubBase + "/c112864f-f6f7-44bb-abc2-93c28e1b4e25.md"
$addrE1 = $githubBase + "/e04279a5-149f-458b-8c7f-6d2d8006abd7.md"

# --- zh-cn sheet ---
Set-HandbackRow $wsZhCn "2" $addrC1 "c112864f-f6f7-44bb-abc2-93c28e1b4e25.md" "c112864f-f6f7-44bb-abc2-93c28e1b4e25.437272b425a9ac945f6d236f5933f008d7c33f30.zh-cn.xlf" "2016-08-18 00:57:57"
Set-HandbackRow $wsZhCn "3" $addrE1 "e04279a5-149f-458b-8c7f-6d2d8006abd7.md" "e04279a5-149f-458b-8c7f-6d2d8006abd7.b0c7ead6a50ae4f6e15fb4428323c53e7d499150.zh-cn.xlf" "2016-08-18 00:57:57"

# --- de-de sheet ---
Set-HandbackRow $wsDeDe "2" $addrC1 "c112864f-f6f7-44bb-abc2-93c28e1b4e25.md" "c112864f-f6f7-44bb-abc2-93c28e1b4e25.437272b425a9ac945f6d236f5933f008d7c33f30.de-de.xlf" "2016-08-18 00:58:10"
Set-HandbackRow $wsDeDe "3" $addrE1 "e04279a5-149f-458b-8c7f-6d2d8006abd7.md" "e04279a5-149f-458b-8c7f-6d2d8006abd7.b0c7ead6a50ae4f6e15fb4428323c53e7d499150.de-de.xlf" "2016-08-18 00:58:10"

# Widen columns C, I, J on both language sheets to fit the newly-populated values.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Range("C1").ColumnWidth = 29.2
    $ws.Range("I1").ColumnWidth = $ws.Range("G1").ColumnWidth
    $ws.Range("J1").ColumnWidth = $ws.Range("G1").ColumnWidth
}

Write-Host "Handback report generated."
